$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values: row3 C becomes "yolo", row4 C becomes "rule", row4 A becomes "test"
# (order chosen so the shared-string table is built as yolo, rule, test)
$ws.Range("C3").Value = "yolo"
$ws.Range("C4").Value = "rule"
$ws.Range("A4").Value = "test"

# Update the selected cell to match the new active selection
$ws.Range("I9").Select()
